$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the subject-count header row (row 1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update meanEMG row (row 2)
$ws.Range("B2").Value = 8.2959479899363604
$ws.Range("C2").Value = 5.9891685282022342
$ws.Range("D2").Value = 7.2116753546531092
$ws.Range("E2").Value = 7.7660648385154882

# Update legmaxROM row (row 3)
$ws.Range("B3").Value = 7.0303179445172486
$ws.Range("C3").Value = 4.8852490163363234
$ws.Range("D3").Value = 7.2281202662674859
$ws.Range("E3").Value = 8.6431981979258197

# Match the updated selection from the edited workbook
$ws.Range("B1:E3").Select()
